$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Find & Replace "U.S.A" with "U./S.A" in the Country column (column C)
# of the financials table - this is what shifts the shared-strings table
# (old "U.S.A" entry removed, new "U./S.A" entry appended at the end).
$rng = $ws.Range("C1:C182")
[void]$rng.Replace("U.S.A", "U./S.A")

# Reflect the final cell selection left by the edit session
[void]$ws.Range("E5").Select()
